# Update the "想去人数" (want-to-go count) column (F) values on the
# "展览" and "全部类型" sheets to reflect the latest generated snapshot.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value = 26
$wsExpo.Range("F5").Value = 8876
$wsExpo.Range("F10").Value = 5407
$wsExpo.Range("F12").Value = 6140
$wsExpo.Range("F15").Value = 403
$wsExpo.Range("F17").Value = 561
$wsExpo.Range("F24").Value = 10007
$wsExpo.Range("F26").Value = 1822
$wsExpo.Range("F29").Value = 2064
$wsExpo.Range("F30").Value = 81
$wsExpo.Range("F34").Value = 2045
$wsExpo.Range("F37").Value = 0
$wsExpo.Range("F47").Value = 1355

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 8876
$wsAll.Range("F13").Value = 5407
$wsAll.Range("F14").Value = 6140
$wsAll.Range("F15").Value = 6140
$wsAll.Range("F18").Value = 403
$wsAll.Range("F20").Value = 561
$wsAll.Range("F28").Value = 10007
$wsAll.Range("F30").Value = 1822
$wsAll.Range("F32").Value = 2064
$wsAll.Range("F33").Value = 81
$wsAll.Range("F48").Value = 1355
